# Fruta / hortaliza, semanal
# Insert a new daily record at the top of the data block (row 181), pushing
# all existing records (rows 181-291) down by one row, and populate the new
# row with the latest reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 181:291 down to 182:292, carrying formatting/styles with them.
$ws.Rows(181).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A181").Value = 4
$ws.Range("B181").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C181").Value = 'Los Lagos'
$ws.Range("D181").Value2 = 44777
$ws.Range("E181").Value = 10
$ws.Range("F181").Value = 100112017
$ws.Range("G181").Value = 'Apio'
$ws.Range("H181").Value = 'Americana (o)'
$ws.Range("I181").Value = 'Primera'
$ws.Range("J181").Value = 35
$ws.Range("K181").Value = 13000
$ws.Range("L181").Value = 14000
$ws.Range("M181").Value = 13571
$ws.Range("N181").Value = '$/docena de matas'
$ws.Range("O181").Value = 'Región de Coquimbo'
$ws.Range("P181").Value = 2262
$ws.Range("Q181").Value = 6
$ws.Range("R181").Value = 'Hortaliza'
